$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "37.032.71"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "2.060.10"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.46"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.672"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.51"
$ws.Range("E8").Value = "  +13.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.04"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.383"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +6.88%  "
$ws.Range("E12").Value = "  +5.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.02"
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").Value = "2.362.64"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("D17").Value = "2.063.61"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "37.017.48"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "0.0₃0931"
$ws.Range("E19").Value = "  +11.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.46"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.24"
$ws.Range("E21").Value = "  +7.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.38"
$ws.Range("E22").Value = "  +2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.21"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.82"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.11"
$ws.Range("E28").Value = "  -5.40%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("E32").Value = "  +7.68%  "
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.42"
$ws.Range("E34").Value = "  +8.02%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.28"
$ws.Range("E37").Value = "  -7.01%  "
$ws.Range("E38").Value = "  -4.86%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +21.83%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.67"
$ws.Range("E42").Value = "  +10.04%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.98"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.15"
$ws.Range("E46").Value = "  +46.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.63"
$ws.Range("E47").Value = "  -52.70%  "
$ws.Range("E48").Value = "  +6.92%  "
$ws.Range("D49").Value = "1.296.22"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.13"
$ws.Range("E51").Value = "  +8.10%  "
